# Finish the excel uploading part: append the newly-uploaded record as row 5
# (and add the function of viewing excel by leaving the selection on the
# freshly written cell D5, mirroring a user who just entered/viewed data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds text-formatted dates (see existing rows which use a "@"
# text number format) - force text formatting before assigning so the new
# value is stored as a shared string like "2014-06-09" instead of being
# auto-converted into a date serial number.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2014-06-09"
$ws.Range("B5").Value = 123
$ws.Range("C5").Value = 18.6
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "hello "

# Reflect the active cell/selection on the newly added row, as captured by
# the workbook's saved view state.
$ws.Range("D5").Select()
